$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.145113
$ws.Range("H2").Value = 6.435339
$ws.Range("Q2").Value = 0.5891066628886666
$ws.Range("R2").Value = 5.301959965998
